{"js": "const newValues = [\n  [\"35+12=\", \"33-5=\", \"84-64=\", \"72-70=\", \"25-18=\"],\n  [\"52-42=\", \"86-7=\", \"96-2=\", \"27+48=\", \"86-2=\"],\n  [\"81-47=\", \"87-33=\", \"11+9=\", \"3+65=\", \"9+22=\"],\n  [\"39+38=\", \"69+17=\", \"58-40=\", \"22+52=\", \"43+7=\"],\n  [\"56+38=\", \"61-14=\", \"55-41=\", \"3+27=\", \"36+14=\"],\n  [\"95-95=\", \"12+7=\", \"26+50=\", \"9+77=\", \"17-6=\"],\n  [\"26+73=\", \"25+35=\", \"78-57=\", \"24+0=\", \"15+79=\"],\n  [\"15+35=\", \"84-22=\", \"0+19=\", \"80-35=\", \"50-3=\"],\n  [\"72+19=\", \"41+14=\", \"12+14=\", \"90-69=\", \"72+14=\"],\n  [\"39-23=\", \"52-49=\", \"14+74=\", \"23+66=\", \"83-33=\"],\n  [\"2+21=\", \"47+11=\", \"90-66=\", \"29+4=\", \"13+24=\"],\n  [\"51+0=\", \"64-18=\", \"20-16=\", \"52-32=\", \"33-20=\"],\n  [\"32-25=\", \"61-28=\", \"7+91=\", \"26+30=\", \"53-26=\"],\n  [\"58-20=\", \"40-34=\", \"89-76=\", \"46+24=\", \"36+61=\"],\n  [\"19+16=\", \"77-10=\", \"6+78=\", \"69-32=\", \"63-29=\"],\n  [\"55+31=\", \"37+28=\", \"4+58=\", \"57-33=\", \"46-36=\"],\n  [\"24+69=\", \"28+33=\", \"97-78=\", \"40-11=\", \"80-29=\"],\n  [\"43+6=\", \"1+93=\", \"97-39=\", \"68-30=\", \"54-43=\"],\n  [\"89-6=\", \"43+39=\", \"83-73=\", \"55+17=\", \"61-20=\"],\n  [\"22-8=\", \"15+45=\", \"10-6=\", \"77+21=\", \"80-8=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(\"Row count mismatch: expected \" + newValues.length + \" got \" + table.rowCount);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$data = @(\n  @('35+12=', '33-5=', '84-64=', '72-70=', '25-18='),\n  @('52-42=', '86-7=', '96-2=', '27+48=', '86-2='),\n  @('81-47=', '87-33=', '11+9=', '3+65=', '9+22='),\n  @('39+38=', '69+17=', '58-40=', '22+52=', '43+7='),\n  @('56+38=', '61-14=', '55-41=', '3+27=', '36+14='),\n  @('95-95=', '12+7=', '26+50=', '9+77=', '17-6='),\n  @('26+73=', '25+35=', '78-57=', '24+0=', '15+79='),\n  @('15+35=', '84-22=', '0+19=', '80-35=', '50-3='),\n  @('72+19=', '41+14=', '12+14=', '90-69=', '72+14='),\n  @('39-23=', '52-49=', '14+74=', '23+66=', '83-33='),\n  @('2+21=', '47+11=', '90-66=', '29+4=', '13+24='),\n  @('51+0=', '64-18=', '20-16=', '52-32=', '33-20='),\n  @('32-25=', '61-28=', '7+91=', '26+30=', '53-26='),\n  @('58-20=', '40-34=', '89-76=', '46+24=', '36+61='),\n  @('19+16=', '77-10=', '6+78=', '69-32=', '63-29='),\n  @('55+31=', '37+28=', '4+58=', '57-33=', '46-36='),\n  @('24+69=', '28+33=', '97-78=', '40-11=', '80-29='),\n  @('43+6=', '1+93=', '97-39=', '68-30=', '54-43='),\n  @('89-6=', '43+39=', '83-73=', '55+17=', '61-20='),\n  @('22-8=', '15+45=', '10-6=', '77+21=', '80-8=')\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\nif ($rowCount -ne $data.Length) {\n    throw \"Row count mismatch: expected \" + $data.Length + \" got \" + $rowCount\n}\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $data[$r-1][$c-1]\n    }\n}\n"}
